$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.822.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.46%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.907.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.60%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5229"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.21%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3786"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.11%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07241"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.12%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9116"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.15%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.58%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07649"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.13%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.916.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.457"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.45%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.00%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008705"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.73%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "27.847.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.49%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.158"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.77%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.193.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.99%  "

# Row 23
$ws.Range("E23").Value = "  +1.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.646"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.56%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.874"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.81%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.170"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.53%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.14%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.867"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09018"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.31%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.865"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.10%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.180"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.67%  "

# Row 34
$ws.Range("E34").Value = "  +1.47%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7810"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.97%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02093"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.81%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.608"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.19%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.076"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.37%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5587"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.34%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.094"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.12%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05287"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "

# Row 42
$ws.Range("E42").Value = "  -2.23%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.44%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.558"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.51%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1515"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4812"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.77%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.19%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.0000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.624"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.73%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.02"
$ws.Range("D50").Style = "Normal"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05994"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
